# Applies the sharedStrings/content updates described by the commit diff.
# Net effect (5 logical cell changes):
#   Metadata!B8   (Date)          : 2025-05-21T14:22:51+00:00 -> 2025-06-13T15:45:04+00:00
#   Metadata!B15  (FHIR Version)  : 4.3.0 -> 4.0.1
#   Elements!AJ2  (Constraint(s) for "Extension")
#       drop the "unless an empty Parameters resource ... or this-is-Parameters" clause
#   Elements!K3   (Type(s) for "Extension.id") : id -> string
#   Elements!M6   (Definition for "Extension.value[x]") : R4B -> R4 link

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

$elements = $wb.Worksheets.Item("Elements")

$elements.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

$elements.Range("K3").Value = "string`n"

$elements.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
